$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4,D5,D6,D7,D8,D9,D10,D11,D12,D14,D15,D16,D17,D19,D20,D23,D24,D25,D26,D27,D28,D29,D30,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51").NumberFormat = "@"

$ws.Range("D2").Value = "31.253.61"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "2.002.37"
$ws.Range("E3").Value = "  +6.50%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "0.7737"
$ws.Range("E5").Value = "  +63.97%  "
$ws.Range("D6").Value = "255.91"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "0.3485"
$ws.Range("E8").Value = "  +20.86%  "
$ws.Range("D9").Value = "28.15"
$ws.Range("E9").Value = "  +27.88%  "
$ws.Range("D10").Value = "0.07000"
$ws.Range("E10").Value = "  +7.20%  "
$ws.Range("D11").Value = "0.8506"
$ws.Range("E11").Value = "  +14.52%  "
$ws.Range("D12").Value = "0.08211"
$ws.Range("E12").Value = "  +5.03%  "
$ws.Range("D13").Value = "1.997.54"
$ws.Range("E13").Value = "  +6.32%  "
$ws.Range("D14").Value = "100.61"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "5.599"
$ws.Range("E15").Value = "  +7.59%  "
$ws.Range("D16").Value = "15.37"
$ws.Range("E16").Value = "  +16.99%  "
$ws.Range("D17").Value = "275.50"
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("D18").Value = "31.249.23"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").Value = "5.907"
$ws.Range("E19").Value = "  +10.19%  "
$ws.Range("D20").Value = "0.000007892"
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("D21").Value = "2.259.12"
$ws.Range("E21").Value = "  +6.77%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "7.110"
$ws.Range("E24").Value = "  +11.46%  "
$ws.Range("D25").Value = "10.02"
$ws.Range("E25").Value = "  +10.07%  "
$ws.Range("D26").Value = "164.74"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").Value = "0.1467"
$ws.Range("E27").Value = "  +51.27%  "
$ws.Range("D28").Value = "19.88"
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("D29").Value = "2.320"
$ws.Range("E29").Value = "  +21.52%  "
$ws.Range("D30").Value = "1.601"
$ws.Range("E30").Value = "  +6.98%  "
$ws.Range("E31").Value = "  +3.26%  "
$ws.Range("D32").Value = "4.608"
$ws.Range("E32").Value = "  +8.37%  "
$ws.Range("D33").Value = "4.413"
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("D34").Value = "0.05211"
$ws.Range("E34").Value = "  +7.87%  "
$ws.Range("D35").Value = "1.232"
$ws.Range("E35").Value = "  +8.85%  "
$ws.Range("D36").Value = "0.7747"
$ws.Range("E36").Value = "  +11.80%  "
$ws.Range("D37").Value = "2.760"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "0.01999"
$ws.Range("E38").Value = "  +5.00%  "
$ws.Range("D39").Value = "2.908"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").Value = "6.679"
$ws.Range("E40").Value = "  +5.40%  "
$ws.Range("D41").Value = "79.34"
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "2.142"
$ws.Range("E42").Value = "  +8.57%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.4671"
$ws.Range("E43").Value = "  +10.11%  "
$ws.Range("D44").Value = "105.96"
$ws.Range("E44").Value = "  +4.57%  "
$ws.Range("D45").Value = "0.8485"
$ws.Range("E45").Value = "  +2.45%  "
$ws.Range("D46").Value = "1.003"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "7.679"
$ws.Range("E47").Value = "  +9.04%  "
$ws.Range("D48").Value = "9.893"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").Value = "0.4309"
$ws.Range("E49").Value = "  +9.26%  "
$ws.Range("D50").Value = "36.77"
$ws.Range("E50").Value = "  +4.80%  "
$ws.Range("D51").Value = "1.516"
$ws.Range("E51").Value = "  +13.35%  "
